$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header years 2009..2020 across B1:M1 ---
$ws.Range("B1").Copy()
$ws.Range("C1:M1").PasteSpecial(-4122)

$ws.Range("B1").Value = 2009
$ws.Range("C1").Value = 2010
$ws.Range("D1").Value = 2011
$ws.Range("E1").Value = 2012
$ws.Range("F1").Value = 2013
$ws.Range("G1").Value = 2014
$ws.Range("H1").Value = 2015
$ws.Range("I1").Value = 2016
$ws.Range("J1").Value = 2017
$ws.Range("K1").Value = 2018
$ws.Range("L1").Value = 2019
$ws.Range("M1").Value = 2020

# --- Row 2: count ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 29
$ws.Range("J2").Value = 27
$ws.Range("K2").Value = 27
$ws.Range("L2").Value = 25
$ws.Range("M2").Value = 26

# --- Row 3: mean ---
$ws.Range("B3:H3").Value = ""
$ws.Range("I3").Value = 1.57677059777363
$ws.Range("J3").Value = 1.853158783125982
$ws.Range("K3").Value = 1.735838798501507
$ws.Range("L3").Value = 1.813199538986507
$ws.Range("M3").Value = 1.853951202552226

# --- Row 4: std ---
$ws.Range("B4:H4").Value = ""
$ws.Range("I4").Value = 1.143308422135791
$ws.Range("J4").Value = 1.097427164630502
$ws.Range("K4").Value = 0.8965988167607843
$ws.Range("L4").Value = 1.373468264939794
$ws.Range("M4").Value = 1.22089507822005

# --- Row 5: min ---
$ws.Range("B5:H5").Value = ""
$ws.Range("I5").Value = 0.6444612049061834
$ws.Range("J5").Value = 0.6897269580981509
$ws.Range("K5").Value = 0.7799170674514385
$ws.Range("L5").Value = 0.6478881242931127
$ws.Range("M5").Value = 0.8066829725159009

# --- Row 6: 25% ---
$ws.Range("B6:H6").Value = ""
$ws.Range("I6").Value = 0.9807633204745193
$ws.Range("J6").Value = 0.9028644183361793
$ws.Range("K6").Value = 0.9190440686592485
$ws.Range("L6").Value = 0.8016320540091433
$ws.Range("M6").Value = 0.9834531181240288

# --- Row 7: 50% ---
$ws.Range("B7:H7").Value = ""
$ws.Range("I7").Value = 1.274529686772878
$ws.Range("J7").Value = 1.844327789636066
$ws.Range("K7").Value = 1.55657204945054
$ws.Range("L7").Value = 1.479470843911362
$ws.Range("M7").Value = 1.167501624428141

# --- Row 8: 75% ---
$ws.Range("B8:H8").Value = ""
$ws.Range("I8").Value = 1.471558570100326
$ws.Range("J8").Value = 2.315881097018852
$ws.Range("K8").Value = 2.282607456583391
$ws.Range("L8").Value = 2.332674655140596
$ws.Range("M8").Value = 2.473009791462784

# --- Row 9: max ---
$ws.Range("B9:H9").Value = ""
$ws.Range("I9").Value = 5.245049700387733
$ws.Range("J9").Value = 4.207530428164726
$ws.Range("K9").Value = 3.769088720112518
$ws.Range("L9").Value = 5.522720497734497
$ws.Range("M9").Value = 4.855903701042041
